$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Delete old row 2 (the "Hiver / Eté / Année" sub-header row).
#    This shifts the data rows (old 3..10) up to become rows 2..9.
# ------------------------------------------------------------------
$ws.Rows(2).Delete()

# ------------------------------------------------------------------
# 2. Build a throw-away named style that reproduces the same font
#    as the existing "s=1" style (Arial 9) but only carries
#    applyFont (no applyNumberFormat). Applying it and then
#    deleting the named style again leaves a clean new cellXfs
#    entry behind that the header cells can reference.
# ------------------------------------------------------------------
$headerStyle = $wb.Styles.Add("TmpHeaderStyle")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9

# ------------------------------------------------------------------
# 3. New header row (row 1).
#    Columns A:E -> plain (default / no special formatting)
#    Columns F:K -> use the new style created above
# ------------------------------------------------------------------
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Remove the temporary named style - the cellXfs entry it left behind
# on F1:K1 remains, but the extra cellStyle/cellStyleXfs bookkeeping
# entries are cleaned up.
$headerStyle.Delete()

# ------------------------------------------------------------------
# 4. Update the sheet selection to match the new layout.
# ------------------------------------------------------------------
$ws.Range("A2:K2").Select()
